$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "added harvard case classification": a new "average_doctor" column is
# inserted immediately before the existing "average_doctor_old" column
# (BP/BQ). The old BP column's data shifts right into BQ, a fresh set of
# values is written into BP, and every "_old" human-doctor-comparison
# series (Ada_old, Avey_old, Babylon_old, Buoy_old, K health_old,
# WebMD_old, doctor_MA_old, doctor_NJ_old, doctor_TH_old) is recomputed
# to reflect the newly added classification case.

$values = @{
  "BP1" = "average_doctor_old"
  "BQ1" = "average_doctor"
  "AI4" = 0.268
  "AJ4" = 0.09
  "AK4" = 0.301
  "AU4" = 0.223
  "AW4" = 0.169
  "BA4" = 1.65
  "BB4" = 0.11
  "BC4" = 0.332
  "BG4" = 0.617
  "BH4" = 0.189
  "BI4" = 0.435
  "BM4" = 0.583
  "BN4" = 0.101
  "BO4" = 0.318
  "BP4" = 0.55
  "BQ4" = 0.678
  "E4" = 0.577
  "F4" = 0.059
  "G4" = 0.243
  "N4" = 0.5
  "O4" = 0.07199999999999999
  "P4" = 0.269
  "Q4" = 0.22
  "R4" = 0.111
  "S4" = 0.333
  "W4" = 0.417
  "X4" = 0.101
  "Y4" = 0.318
  "AI5" = 0.282
  "AJ5" = 0.09
  "AK5" = 0.3
  "AU5" = 0.385
  "AV5" = 0.105
  "AW5" = 0.324
  "BA5" = 1.029
  "BB5" = 0.031
  "BC5" = 0.177
  "BG5" = 0.377
  "BH5" = 0.06900000000000001
  "BI5" = 0.263
  "BM5" = 0.337
  "BN5" = 0.031
  "BO5" = 0.176
  "BP5" = 0.343
  "BQ5" = 0.383
  "E5" = 0.615
  "F5" = 0.055
  "G5" = 0.235
  "N5" = 0.718
  "O5" = 0.046
  "P5" = 0.214
  "Q5" = 0.145
  "R5" = 0.045
  "S5" = 0.212
  "W5" = 0.295
  "X5" = 0.082
  "Y5" = 0.286
  "AI6" = 0.275
  "AU6" = 0.282
  "BA6" = 1.266
  "BG6" = 0.468
  "BM6" = 0.427
  "BP6" = 0.422
  "BQ6" = 0.488
  "E6" = 0.595
  "N6" = 0.589
  "Q6" = 0.175
  "W6" = 0.346
  "AI7" = 0.279
  "AU7" = 0.336
  "BA7" = 1.112
  "BG7" = 0.409
  "BM7" = 0.368
  "BP7" = 0.371
  "BQ7" = 0.419
  "E7" = 0.607
  "N7" = 0.66
  "Q7" = 0.156
  "W7" = 0.313
  "AI8" = 0.292
  "AJ8" = 0.131
  "AK8" = 0.362
  "AU8" = 0.343
  "AV8" = 0.108
  "AW8" = 0.328
  "BA8" = 1.525
  "BB8" = 0.102
  "BC8" = 0.32
  "BG8" = 0.507
  "BH8" = 0.144
  "BI8" = 0.379
  "BM8" = 0.591
  "BN8" = 0.078
  "BO8" = 0.278
  "BP8" = 0.508
  "BQ8" = 0.5610000000000001
  "E8" = 0.648
  "F8" = 0.08
  "G8" = 0.283
  "N8" = 0.901
  "O8" = 0.008
  "P8" = 0.089
  "Q8" = 0.159
  "R8" = 0.08400000000000001
  "S8" = 0.29
  "W8" = 0.477
  "X8" = 0.111
  "Y8" = 0.333
  "AI9" = 0.3
  "AJ9" = 0.21
  "AK9" = 0.458
  "BA9" = 1.6
  "BM9" = 0.7
  "BN9" = 0.21
  "BO9" = 0.458
  "BP9" = 0.533
  "BQ9" = 0.581
  "E9" = 0.5
  "F9" = 0.25
  "G9" = 0.5
  "N9" = 1
  "O9" = 0
  "P9" = 0
  "AI10" = 0.3
  "AJ10" = 0.21
  "AK10" = 0.458
  "BA10" = 1.9
  "BB10" = 0.25
  "BC10" = 0.5
  "BM10" = 0.8
  "BN10" = 0.16
  "BO10" = 0.4
  "BP10" = 0.633
  "BQ10" = 0.714
  "E10" = 0.6
  "F10" = 0.24
  "G10" = 0.49
  "N10" = 1
  "O10" = 0
  "P10" = 0
  "W10" = 0.6
  "X10" = 0.24
  "Y10" = 0.49
  "AI11" = 0.3
  "AJ11" = 0.21
  "AK11" = 0.458
  "AU11" = 0.4
  "AV11" = 0.24
  "AW11" = 0.49
  "BA11" = 1.9
  "BB11" = 0.25
  "BC11" = 0.5
  "BM11" = 0.8
  "BN11" = 0.16
  "BO11" = 0.4
  "BP11" = 0.633
  "BQ11" = 0.714
  "E11" = 0.7
  "F11" = 0.21
  "G11" = 0.458
  "N11" = 1
  "O11" = 0
  "P11" = 0
  "W11" = 0.6
  "X11" = 0.24
  "Y11" = 0.49
  "AU12" = 3.4
  "AV12" = 5.04
  "AW12" = 2.245
  "BA12" = 3.658
  "BB12" = 0.16
  "BC12" = 0.4
  "BM12" = 1.125
  "BN12" = 0.109
  "BO12" = 0.331
  "BP12" = 1.219
  "BQ12" = 1.261
  "E12" = 1.857
  "F12" = 2.122
  "G12" = 1.457
  "W12" = 1.167
  "X12" = 0.139
  "Y12" = 0.373
  "BP13" = 0.714
  "BQ13" = 0.643
}

foreach ($cellRef in $values.Keys) {
  $ws.Range($cellRef).Value = $values[$cellRef]
}
